# Update "F" column (想去人数 / "people interested") values on the
# "展览" and "全部类型" worksheets to match the newly generated data.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 1867
$ws1.Range("F6").Value  = 2567
$ws1.Range("F7").Value  = 169
$ws1.Range("F10").Value = 1524
$ws1.Range("F11").Value = 528
$ws1.Range("F12").Value = 43
$ws1.Range("F13").Value = 330
$ws1.Range("F14").Value = 230
$ws1.Range("F17").Value = 209
$ws1.Range("F20").Value = 13
$ws1.Range("F21").Value = 173
$ws1.Range("F22").Value = 56
$ws1.Range("F23").Value = 1637
$ws1.Range("F25").Value = 400
$ws1.Range("F26").Value = 570
$ws1.Range("F27").Value = 204
$ws1.Range("F28").Value = 297
$ws1.Range("F29").Value = 416

# --- Sheet "全部类型" (all types, contains an extra row vs. 展览) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 1867
$ws4.Range("F7").Value  = 2567
$ws4.Range("F8").Value  = 169
$ws4.Range("F11").Value = 1524
$ws4.Range("F12").Value = 528
$ws4.Range("F13").Value = 43
$ws4.Range("F14").Value = 330
$ws4.Range("F15").Value = 230
$ws4.Range("F18").Value = 209
$ws4.Range("F21").Value = 13
$ws4.Range("F22").Value = 173
$ws4.Range("F23").Value = 56
$ws4.Range("F24").Value = 1637
$ws4.Range("F26").Value = 400
$ws4.Range("F27").Value = 570
$ws4.Range("F28").Value = 204
$ws4.Range("F29").Value = 297
$ws4.Range("F30").Value = 416
